# Adding support for multi slide bhajans
#
# Slide 1 carries a single bhajan's title + lyrics across two textboxes:
#   "Shape 55" - title line, originally two runs:
#                "Bhajo Mana Krishna " + "Gopal"
#                -> collapses to a single run reading "Bhajo"
#                   (the lyric is now continued on a later/other slide,
#                   so the title is shortened to just the first word).
#   "Shape 56" - lyrics, originally three paragraphs:
#                "O Mind! Chant the name of Krishna; who is known as
#                 Gopala, Mukunda and Govinda; " / "beloved son of Nanda"
#                 / (empty)
#                -> collapses to a single paragraph reading "O Mind! "
#                   (the rest of the lyric moves to the continuation
#                   slide(s) that make multi-slide bhajans possible).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Find-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $null
}

# --- "Shape 55": "Bhajo Mana Krishna " + "Gopal" -> "Bhajo" ------------
$shape55 = Find-ShapeByName $s "Shape 55"
$titleRange = $shape55.TextFrame.TextRange

# Remove the first run ("Bhajo Mana Krishna ") in its entirety, leaving
# the second run ("Gopal") - along with its own distinct run formatting
# (it carries smtClean="0") - as the sole surviving run in the paragraph.
$splitAt = $titleRange.Text.IndexOf("Gopal")
$firstRun = $titleRange.Characters(1, $splitAt)
$firstRun.Text = ""

# Re-fetch the (now shortened) range and rename what is left ("Gopal")
# to "Bhajo", keeping that run's own formatting untouched.
$titleRange = $shape55.TextFrame.TextRange
$remainingRun = $titleRange.Characters(1, $titleRange.Length)
$remainingRun.Text = "Bhajo"

# --- "Shape 56": trim the lyric down to its first sentence -------------
$shape56 = Find-ShapeByName $s "Shape 56"
$lyricRange = $shape56.TextFrame.TextRange

# Keep only "O Mind! " and drop everything else (rest of paragraph 1,
# all of paragraph 2, and the empty trailing paragraph). The deleted
# trailing/empty paragraph's endParaRPr (sz=2800) formatting is carried
# over onto the now-merged single paragraph.
$keepText = "O Mind! "
$tail = $lyricRange.Characters($keepText.Length + 1, $lyricRange.Length)
$tail.Delete()
